$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.201.30'
$ws.Range("E2").Value = '  +0.73%  '
$ws.Range("D3").Value = '1.562.09'
$ws.Range("D4").Value = '''1.01'
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = '''210.24'
$ws.Range("E5").Value = '  +1.26%  '
$ws.Range("E6").Value = '  +0.13%  '
$ws.Range("D7").Value = '''1.00'
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("D8").Value = '''21.99'
$ws.Range("E8").Value = '  -0.71%  '
$ws.Range("E9").Value = '  -0.05%  '
$ws.Range("D10").Value = '''0.0596'
$ws.Range("E10").Value = '  -1.07%  '
$ws.Range("D11").Value = '''0.0872'
$ws.Range("E11").Value = '  +1.62%  '
$ws.Range("D12").Value = '1.784.75'
$ws.Range("E12").Value = '  -0.14%  '
$ws.Range("D13").Value = '1.567.79'
$ws.Range("E13").Value = '  +0.19%  '
$ws.Range("E14").Value = '  +0.32%  '
$ws.Range("E15").Value = '  -0.58%  '
$ws.Range("D16").Value = '27.171.14'
$ws.Range("E16").Value = '  +0.62%  '
$ws.Range("D17").Value = '''61.83'
$ws.Range("E17").Value = '  -0.30%  '
$ws.Range("B18").Value = 'Chainlink'
$ws.Range("C18").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D18").Value = '''7.43'
$ws.Range("E18").Value = '  +1.01%  '
$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").Value = '''216.32'
$ws.Range("E19").Value = '  +0.05%  '
$ws.Range("D20").Value = '0.0₃0701'
$ws.Range("E20").Value = '  -0.54%  '
$ws.Range("E21").Value = '  -0.19%  '
$ws.Range("E22").Value = '  +0.24%  '
$ws.Range("D23").Value = '''9.21'
$ws.Range("E23").Value = '  -0.05%  '
$ws.Range("D25").Value = '''153.14'
$ws.Range("E25").Value = '  +0.00%  '
$ws.Range("D26").Value = '''6.63'
$ws.Range("E26").Value = '  +0.05%  '
$ws.Range("E27").Value = '  -0.79%  '
$ws.Range("E28").Value = '  +1.45%  '
$ws.Range("E29").Value = '  -0.20%  '
$ws.Range("E30").Value = '  +2.19%  '
$ws.Range("D31").Value = '''0.0470'
$ws.Range("E31").Value = '  -0.10%  '
$ws.Range("E32").Value = '  -0.14%  '
$ws.Range("E33").Value = '  +1.05%  '
$ws.Range("D34").Value = '1.434.15'
$ws.Range("E34").Value = '  +0.79%  '
$ws.Range("E35").Value = '  +2.80%  '
$ws.Range("E36").Value = '  -0.06%  '
$ws.Range("D37").Value = '''2.34'
$ws.Range("E37").Value = '  -0.46%  '
$ws.Range("E38").Value = '  +0.47%  '
$ws.Range("E39").Value = '  -0.12%  '
$ws.Range("D40").Value = '''5.92'
$ws.Range("E40").Value = '  +1.73%  '
$ws.Range("D41").Value = '''0.806'
$ws.Range("E41").Value = '  -0.28%  '
$ws.Range("E42").Value = '  -0.14%  '
$ws.Range("E43").Value = '  +0.13%  '
$ws.Range("D44").Value = '''0.998'
$ws.Range("E44").Value = '  -0.81%  '
$ws.Range("D45").Value = '''64.29'
$ws.Range("E45").Value = '  -0.72%  '
$ws.Range("D46").Value = '''1.73'
$ws.Range("E46").Value = '  -1.35%  '
$ws.Range("D47").Value = '1.698.89'
$ws.Range("E47").Value = '  -0.06%  '
$ws.Range("D48").Value = '''85.65'
$ws.Range("E48").Value = '  -2.01%  '
$ws.Range("D49").Value = '''0.0524'
$ws.Range("E49").Value = '  +0.86%  '
$ws.Range("D50").Value = '0.0₆0102'
$ws.Range("E50").Value = '  +1.92%  '
$ws.Range("D51").Value = '''0.0949'
$ws.Range("E51").Value = '  -1.10%  '
